$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the total "Valor Mora" summary cell (E11): 341509 -> 29509
$ws.Range("E11").Value2 = 29509

# Update "Cant. Trabajadores" (C13): 2 -> 1
$ws.Range("C13").Value2 = 1

# Update "Cant. Periodos" (F13): 7 -> 1
$ws.Range("F13").Value2 = 1

# Update the remaining worker's "Salario Basico" value (G16): 737717 -> 781242
$ws.Range("G16").Value2 = 781242

# Remove the second worker's block (IRINA PAOLA INSIGNARES OYOLA, periods 2502-2507)
# which occupied rows 17-22, shifting the signature block up.
$ws.Range("B17:J22").EntireRow.Delete()
